# draft-gandhi-mpls-ioam-sr-05.pptx -- "Add files via upload" edit replay
#
# Changes applied:
#  - Slide 13 ("HbH Indicator Label Allocation Methods"): table header cell
#    "Extra Label stack Size" -> "Extra Label stack Size (Note 4)"
#  - Slide 9  ("E2E Indicator Label Allocation Methods"): same table header
#    cell -> "Extra Label stack Size (Note 2)"
#  - Slide 14: widen/reposition the Content Placeholder to span the full
#    slide width (off/ext change only, text untouched)
#  - Slide 19 ("4. IOAM Data After EOS"):
#      * re-center/resize the title box and drop the leading "4. " from the
#        title text
#      * tidy trailing punctuation/spacing on three bullet paragraphs

$p = $ppt.ActivePresentation

# --- Slide 13: table cell text -------------------------------------------
$s13 = $p.Slides.Item(13)
$tbl13 = $s13.Shapes.Item(2).Table
$tbl13.Cell(1, 2).Shape.TextFrame.TextRange.Text = "Extra Label stack Size (Note 4)"

# --- Slide 14: Content Placeholder resize/reposition ----------------------
$s14 = $p.Slides.Item(14)
$contentPh14 = $s14.Shapes.Item(2)
$contentPh14.Left = 36
$contentPh14.Width = 648

# --- Slide 19: title box resize + text, and bullet tidy-up ----------------
$s19 = $p.Slides.Item(19)

$title19 = $s19.Shapes.Item(1)
$title19.Left = 168
$title19.Width = 288
$title19.TextFrame.TextRange.Text = "IOAM Data After EOS"

$body19 = $s19.Shapes.Item(2)
$tr19 = $body19.TextFrame.TextRange

$tr19.Paragraphs(4, 1).Runs(1, 1).Text = "In all these cases there is a CW immediately after EOS "
$tr19.Paragraphs(5, 1).Runs(1, 1).Text = "Then there is the universal fragmentation idea that is floating about that also wants to follow EOS"
$tr19.Paragraphs(8, 1).Runs(1, 1).Text = "This is a generic issue applicable to all G-ACH mechanisms used for data traffic"

# --- Slide 9: table cell text ---------------------------------------------
$s9 = $p.Slides.Item(9)
$tbl9 = $s9.Shapes.Item(2).Table
$tbl9.Cell(1, 2).Shape.TextFrame.TextRange.Text = "Extra Label stack Size (Note 2)"
